$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.05297061800956726
$ws.Range("B2").Value = 0.9850805997848511
$ws.Range("C2").Value = 0.02776168286800385
$ws.Range("D2").Value = 0.9953567981719971

$ws.Range("A3").Value = 0.008638261817395687
$ws.Range("B3").Value = 0.9984014630317688
$ws.Range("C3").Value = 0.0158476997166872
$ws.Range("D3").Value = 0.9963343143463135

$ws.Range("A4").Value = 0.004096722695976496
$ws.Range("B4").Value = 0.9988391995429993
$ws.Range("C4").Value = 0.002228233031928539
$ws.Range("D4").Value = 0.9990224838256836

$ws.Range("A5").Value = 0.001423759269528091
$ws.Range("B5").Value = 0.9996384382247925
$ws.Range("C5").Value = 0.0006458171410486102
$ws.Range("D5").Value = 1

$ws.Range("A6").Value = 0.001939180307090282
$ws.Range("B6").Value = 0.9996194243431091
$ws.Range("C6").Value = 0.000266658200416714
$ws.Range("D6").Value = 1

$ws.Range("A7").Value = 0.001024800701998174
$ws.Range("B7").Value = 0.9997145533561707
$ws.Range("C7").Value = 0.00008454782073386014
$ws.Range("D7").Value = 1

$ws.Range("A8").Value = 0.00196632812730968
$ws.Range("B8").Value = 0.9995813369750977
$ws.Range("C8").Value = 0.0001295774127356708
$ws.Range("D8").Value = 1

$ws.Range("A9").Value = 0.0007570366724394262
$ws.Range("B9").Value = 0.9997716546058655
$ws.Range("C9").Value = 0.0000506466931256
$ws.Range("D9").Value = 1

$ws.Range("A10").Value = 0.0007868724060244858
$ws.Range("B10").Value = 0.9997716546058655
$ws.Range("C10").Value = 0.00004302233355701901
$ws.Range("D10").Value = 1

$ws.Range("A11").Value = 0.0005520881386473775
$ws.Range("B11").Value = 0.9998857975006104
$ws.Range("C11").Value = 0.00001558025724079926
$ws.Range("D11").Value = 1

$ws.Range("A12").Value = 0.0005973918596282601
$ws.Range("B12").Value = 0.9998287558555603
$ws.Range("C12").Value = 0.00003197484329575673
$ws.Range("D12").Value = 1

$ws.Range("A13").Value = 0.0002919227408710867
$ws.Range("B13").Value = 0.9999238848686218
$ws.Range("C13").Value = 0.000005278091066429624
$ws.Range("D13").Value = 1

$ws.Range("A14").Value = 0.0005308371037244797
$ws.Range("B14").Value = 0.9998096823692322
$ws.Range("C14").Value = 0.000007331217148021096
$ws.Range("D14").Value = 1

$ws.Range("A15").Value = 0.0006463845493271947
$ws.Range("B15").Value = 0.9999048709869385
$ws.Range("C15").Value = 0.00009279639925807714
$ws.Range("D15").Value = 1

$ws.Range("A16").Value = 0.0006180675700306892
$ws.Range("B16").Value = 0.999866783618927
$ws.Range("C16").Value = 0.0000009865923402685439
$ws.Range("D16").Value = 1

$ws.Range("A17").Value = 0.0005904207355342805
$ws.Range("B17").Value = 0.9998477697372437
$ws.Range("C17").Value = 0.000001955155084942817
$ws.Range("D17").Value = 1

$ws.Range("A18").Value = 0.0002343400556128472
$ws.Range("B18").Value = 0.9999619126319885
$ws.Range("C18").Value = 0.000001502354052718147
$ws.Range("D18").Value = 1

$ws.Range("A19").Value = 0.0006626341491937637
$ws.Range("B19").Value = 0.9998096823692322
$ws.Range("C19").Value = 0.0000004290456843136781
$ws.Range("D19").Value = 1

$ws.Range("A20").Value = 0.0007167106959968805
$ws.Range("B20").Value = 0.9998857975006104
$ws.Range("C20").Value = 0.0000005285048132463999
$ws.Range("D20").Value = 1

$ws.Range("A21").Value = 0.0001745160552673042
$ws.Range("B21").Value = 0.9999619126319885
$ws.Range("C21").Value = 0.0000003326848911910929
$ws.Range("D21").Value = 1

$ws.Range("A22").Value = 0.0001081222653738223
$ws.Range("B22").Value = 0.9999619126319885
$ws.Range("C22").Value = 0.0000001158743714313459
$ws.Range("D22").Value = 1

$ws.Range("A23").Value = 0.000142380129545927
$ws.Range("B23").Value = 0.9999809861183167
$ws.Range("C23").Value = 0.0000000832798221495068
$ws.Range("D23").Value = 1

$ws.Range("A24").Value = 0.0007438646862283349
$ws.Range("B24").Value = 0.9998096823692322
$ws.Range("C24").Value = 0.0000007302297149180959
$ws.Range("D24").Value = 1

$ws.Range("A25").Value = 0.0001643219875404611
$ws.Range("B25").Value = 0.9999238848686218
$ws.Range("C25").Value = 0.0000001944291057043301
$ws.Range("D25").Value = 1

$ws.Range("A26").Value = 0.0000855125836096704
$ws.Range("B26").Value = 0.9999619126319885
$ws.Range("C26").Value = 0.00000002336387794343864
$ws.Range("D26").Value = 1

$ws.Range("A27").Value = 0.00000280625090454123
$ws.Range("B27").Value = 1
$ws.Range("C27").Value = 0.00000000949708045538955
$ws.Range("D27").Value = 1

$ws.Range("A28").Value = 0.00009868038614513353
$ws.Range("B28").Value = 0.9999428987503052
$ws.Range("C28").Value = 0.0000000140416824834233
$ws.Range("D28").Value = 1

$ws.Range("A29").Value = 0.0001325308840023354
$ws.Range("B29").Value = 0.9999619126319885
$ws.Range("C29").Value = 0.00000004180414947541067
$ws.Range("D29").Value = 1

$ws.Range("A30").Value = 0.0005380103830248117
$ws.Range("B30").Value = 0.9998857975006104
$ws.Range("C30").Value = 0.00000007751712871595373
$ws.Range("D30").Value = 1

$ws.Range("A31").Value = 0.0006957969162613153
$ws.Range("B31").Value = 0.9998477697372437
$ws.Range("C31").Value = 0.00000002391717046634767
$ws.Range("D31").Value = 1

$ws.Range("A32").Value = 0.0001476978650316596
$ws.Range("B32").Value = 0.9999619126319885
$ws.Range("C32").Value = 0.00000002144094679579212
$ws.Range("D32").Value = 1

$ws.Range("A33").Value = 0.0004551385645754635
$ws.Range("B33").Value = 0.9999428987503052
$ws.Range("C33").Value = 0.00000008415681662654606
$ws.Range("D33").Value = 1

$ws.Range("A34").Value = 0.00002671026231837459
$ws.Range("B34").Value = 0.9999809861183167
$ws.Range("C34").Value = 0.00000001890647816082947
$ws.Range("D34").Value = 1

$ws.Range("A35").Value = 0.0000256000566878356
$ws.Range("B35").Value = 0.9999809861183167
$ws.Range("C35").Value = 0.00000002414996913557843
$ws.Range("D35").Value = 1

$ws.Range("A36").Value = 0.00004637566962628625
$ws.Range("B36").Value = 0.9999619126319885
$ws.Range("C36").Value = 0.000000009176627457918585
$ws.Range("D36").Value = 1

$ws.Range("A37").Value = 0.0004940321668982506
$ws.Range("B37").Value = 0.9999238848686218
$ws.Range("C37").Value = 0.0000007920424422991346
$ws.Range("D37").Value = 1

$ws.Range("A38").Value = 0.0001078139684977941
$ws.Range("B38").Value = 0.9999428987503052
$ws.Range("C38").Value = 0.000000003350197230744811
$ws.Range("D38").Value = 1

$ws.Range("A39").Value = 0.0002516054082661867
$ws.Range("B39").Value = 0.9999619126319885
$ws.Range("C39").Value = 0.000002085740561597049
$ws.Range("D39").Value = 1

$ws.Range("A40").Value = 0.0006333301425911486
$ws.Range("B40").Value = 0.9998477697372437
$ws.Range("C40").Value = 0.00000003315205532317123
$ws.Range("D40").Value = 1

$ws.Range("A41").Value = 0.0001859597105067223
$ws.Range("B41").Value = 0.9999048709869385
$ws.Range("C41").Value = 0.00000002426692802259822
$ws.Range("D41").Value = 1

$ws.Range("A42").Value = 0.000003200391120117274
$ws.Range("B42").Value = 1
$ws.Range("C42").Value = 0.00000002846186752947233
$ws.Range("D42").Value = 1

$ws.Range("A43").Value = 0.000005458916348288767
$ws.Range("B43").Value = 1
$ws.Range("C43").Value = 0.000000007894823461640499
$ws.Range("D43").Value = 1

$ws.Range("A44").Value = 0.0002096290118061006
$ws.Range("B44").Value = 0.9999428987503052
$ws.Range("C44").Value = 0.00000014451138952154
$ws.Range("D44").Value = 1

$ws.Range("A45").Value = 0.0004857924941461533
$ws.Range("B45").Value = 0.9998857975006104
$ws.Range("C45").Value = 0.0000003613719457007392
$ws.Range("D45").Value = 1

$ws.Range("A46").Value = 0.0001506252883700654
$ws.Range("B46").Value = 0.9999619126319885
$ws.Range("C46").Value = 0.00001468274331273278
$ws.Range("D46").Value = 1

$ws.Range("A47").Value = 0.0003001784207299352
$ws.Range("B47").Value = 0.9999428987503052
$ws.Range("C47").Value = 0.00008874215563992038
$ws.Range("D47").Value = 1

$ws.Range("A48").Value = 0.0002384933468420058
$ws.Range("B48").Value = 0.9999428987503052
$ws.Range("C48").Value = 0.00000006988642553551472
$ws.Range("D48").Value = 1

$ws.Range("A49").Value = 0.00008121335849864408
$ws.Range("B49").Value = 0.9999809861183167
$ws.Range("C49").Value = 0.00000001377949754299834
$ws.Range("D49").Value = 1

$ws.Range("A50").Value = 0.0000565636764804367
$ws.Range("B50").Value = 0.9999809861183167
$ws.Range("C50").Value = 0.00000001223543488748646
$ws.Range("D50").Value = 1

$ws.Range("A51").Value = 0.0001373050617985427
$ws.Range("B51").Value = 0.9999619126319885
$ws.Range("C51").Value = 0.0000002176775524276309
$ws.Range("D51").Value = 1

